$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded data is the same grid scaled up by a factor of 7.6
# (e.g. row/column totals were re-derived from the same base shape).
for ($r = 1; $r -le 9; $r++) {
    for ($c = 1; $c -le 15; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value2 = $cell.Value2 * 7.6
    }
}

# Reflect the saved selection state: row 1 through row 23 selected,
# with the active cell on row 23.
[void]$ws.Range("A1:XFD23").Select()
